$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered list of values for A2:A28 (row 29 / value 69 removed,
# and 3 / 40 rotated to the end of the list).
$values = @(51,81,143,145,160,223,239,277,307,324,353,513,529,568,765,1216,1222,1228,1234,1325,1332,1455,1533,489,1357,3,40)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Remove the now-unused last row (previously row 29).
$ws.Rows.Item(29).Delete()
